$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:J1) text: drop the " Seen Rx" suffix / fix casing ---
$headerValues = @("FFTR", "LIGAZID", "EMAZID", "LIPICON", "AGLIP", "CIFIBET", "AMLEVO", "CARDOBIS", "RIVAROX", "NOCLOG")
for ($i = 0; $i -lt $headerValues.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerValues[$i]
}

# --- Full data block (A2:J41): label in col A, 9 numeric values in B:J ---
# Rows are now sorted by CCF label and include 9 new CCF rows (CCF, CCF15, CCF25, CCF36, CCF43, CCF52, CCF55, CCF56, CCF58)
$data = @(
    @("CCF", 39, 29, 11, 5, 6, 25, 56, 8, 41),
    @("CCF10", 10, 8, 3, 4, 2, 4, 17, 2, 13),
    @("CCF11", 0, 1, 0, 0, 0, 1, 1, 0, 1),
    @("CCF12", 2, 2, 2, 0, 0, 0, 1, 1, 2),
    @("CCF13", 2, 2, 0, 2, 1, 1, 5, 0, 3),
    @("CCF14", 0, 0, 0, 0, 0, 0, 1, 0, 0),
    @("CCF15", 2, 1, 0, 0, 0, 0, 3, 0, 0),
    @("CCF16", 2, 2, 0, 2, 1, 2, 4, 0, 6),
    @("CCF17", 2, 0, 1, 0, 0, 0, 2, 1, 1),
    @("CCF20", 8, 7, 6, 0, 2, 4, 17, 2, 10),
    @("CCF21", 0, 0, 1, 0, 0, 1, 1, 0, 0),
    @("CCF22", 3, 6, 3, 0, 0, 0, 11, 0, 4),
    @("CCF23", 4, 1, 0, 0, 0, 0, 0, 0, 0),
    @("CCF24", 1, 0, 2, 0, 1, 0, 3, 0, 2),
    @("CCF25", 0, 0, 0, 0, 0, 2, 2, 2, 3),
    @("CCF26", 0, 0, 0, 0, 1, 1, 0, 0, 1),
    @("CCF30", 2, 5, 0, 0, 0, 3, 8, 1, 1),
    @("CCF31", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF32", 1, 0, 0, 0, 0, 1, 3, 1, 1),
    @("CCF33", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF34", 1, 5, 0, 0, 0, 0, 4, 0, 0),
    @("CCF35", 0, 0, 0, 0, 0, 2, 1, 0, 0),
    @("CCF36", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF40", 17, 9, 1, 1, 1, 13, 12, 2, 16),
    @("CCF41", 1, 0, 1, 0, 0, 1, 4, 0, 14),
    @("CCF42", 10, 9, 0, 0, 0, 1, 1, 0, 0),
    @("CCF43", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF44", 1, 0, 0, 0, 1, 3, 1, 0, 0),
    @("CCF45", 0, 0, 0, 0, 0, 1, 5, 0, 0),
    @("CCF46", 5, 0, 0, 0, 0, 7, 1, 2, 2),
    @("CCF47", 0, 0, 0, 1, 0, 0, 0, 0, 0),
    @("CCF50", 2, 0, 1, 0, 1, 1, 2, 1, 1),
    @("CCF51", 0, 0, 1, 0, 0, 0, 0, 0, 1),
    @("CCF52", 1, 0, 0, 0, 0, 0, 1, 1, 0),
    @("CCF53", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF54", 0, 0, 0, 0, 1, 0, 0, 0, 0),
    @("CCF55", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF56", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("CCF57", 1, 0, 0, 0, 0, 1, 1, 0, 0),
    @("CCF58", 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

$startRow = 2
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowArr = $data[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $rowArr.Count; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowArr[$c]
    }
}

Write-Output "done"
